$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C10").Value = "new Integer(p.getPartId()).toString()"
$chars = $ws4.Range("C10").Characters(5, 34)
$chars.Font.Color = 13023145
$chars.Font.Name = "JetBrains Mono"
Write-Output "done"
